$d = $word.ActiveDocument

# --- Step 1: Update the "This sample is compatible with the ..." paragraph ---
# Merge the two runs into one, update the SDK version text, and apply the new
# paragraph/run formatting (Heading1 style, spacing before=0, and rFonts/color/
# sz/szCs additions on top of the existing italic).
$compatPara = $d.Paragraphs(2)
$markRange = $d.Range($compatPara.Range.End - 1, $compatPara.Range.End)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
  '<w:body>' + `
  '<w:p w14:paraId="5D973B3C" w14:textId="77777777" w:rsidR="000A7493" w:rsidRDefault="007B34BA" w:rsidP="000A169A">' + `
  '<w:pPr>' + `
  '<w:pStyle w:val="Heading1"/>' + `
  '<w:spacing w:before="0"/>' + `
  '<w:rPr>' + `
  '<w:rFonts w:eastAsiaTheme="minorHAnsi" w:cs="Times New Roman"/>' + `
  '<w:i/>' + `
  '<w:color w:val="auto"/>' + `
  '<w:sz w:val="20"/>' + `
  '<w:szCs w:val="22"/>' + `
  '</w:rPr>' + `
  '</w:pPr>' + `
  '<w:r>' + `
  '<w:rPr>' + `
  '<w:rFonts w:eastAsiaTheme="minorHAnsi" w:cs="Times New Roman"/>' + `
  '<w:i/>' + `
  '<w:color w:val="auto"/>' + `
  '<w:sz w:val="20"/>' + `
  '<w:szCs w:val="22"/>' + `
  '</w:rPr>' + `
  '<w:t>This sample is compatible with the Windows 10 April 2018 Update SDK (17134)</w:t>' + `
  '</w:r>' + `
  '</w:p>' + `
  '</w:body>' + `
  '</w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'
$markRange.InsertXML($xml)

# --- Step 2: Move the "_GoBack" bookmark to the now-empty paragraph that
# follows (it used to sit on the "This sample is set up to require..."
# paragraph further down; Word re-stamps it at the last edit location, which
# is this blank paragraph right after the text we just changed). Adding a
# bookmark with a name that already exists elsewhere moves it here and
# removes it from its old location. ---
$blankPara = $d.Paragraphs(3)
$d.Bookmarks.Add("_GoBack", $blankPara.Range)
